$wb = $excel.ActiveWorkbook

# --- Sheet: Excel_vs_ML ---
$ws = $wb.Worksheets.Item("Excel_vs_ML")

$ws.Range("H2").Value = 149115.72
$ws.Range("L2").Value = 148.6
$ws.Range("M2").Value = -48765.35
$ws.Range("H3").Value = 514393.08
$ws.Range("L3").Value = 113.32
$ws.Range("M3").Value = 38201.55
$ws.Range("N3").Value = 3820.15
$ws.Range("S3").Value = 38201.54999999999
$ws.Range("H4").Value = 408711.91
$ws.Range("L4").Value = 138.76
$ws.Range("M4").Value = -114171.53
$ws.Range("H5").Value = 760006.08
$ws.Range("L5").Value = 149.29
$ws.Range("M5").Value = -250908.71
$ws.Range("H6").Value = 362804.55
$ws.Range("L6").Value = 94.87
$ws.Range("M6").Value = 19604.7
$ws.Range("H7").Value = 528535.6899999999
$ws.Range("L7").Value = 100.54
$ws.Range("M7").Value = -2839.75
$ws.Range("O7").Value = "On Track"
$ws.Range("H8").Value = 38505.81
$ws.Range("L8").Value = 92.23
$ws.Range("M8").Value = 98335.11
$ws.Range("N8").Value = 2398.42
$ws.Range("S8").Value = 98335.11000000002
$ws.Range("H10").Value = 115615.11
$ws.Range("L10").Value = 97.88
$ws.Range("M10").Value = 333243.04
$ws.Range("N10").Value = 7934.36
$ws.Range("O10").Value = "On Track"
$ws.Range("S10").Value = 333243.04
$ws.Range("T10").Value = "YES"
$ws.Range("H12").Value = 285615.72
$ws.Range("L12").Value = 110.1
$ws.Range("M12").Value = 274191.01
$ws.Range("N12").Value = 6231.61
$ws.Range("O12").Value = "Overpacing"
$ws.Range("S12").Value = 274191.01
$ws.Range("T12").Value = "NO"
$ws.Range("H13").Value = 47892.43
$ws.Range("L13").Value = 95.43000000000001
$ws.Range("M13").Value = 121485.11
$ws.Range("N13").Value = 2131.32
$ws.Range("O13").Value = "On Track"
$ws.Range("S13").Value = 121485.11
$ws.Range("T13").Value = "YES"
$ws.Range("H14").Value = 322426.62
$ws.Range("L14").Value = 96.27
$ws.Range("M14").Value = 12509.47
$ws.Range("O14").Value = "On Track"
$ws.Range("H16").Value = 195465.22
$ws.Range("L16").Value = 120.95
$ws.Range("M16").Value = -33860.12
$ws.Range("H17").Value = 123864.18
$ws.Range("L17").Value = 129.59
$ws.Range("M17").Value = -28280.37
$ws.Range("H20").Value = 221276.88
$ws.Range("L20").Value = 120.52
$ws.Range("M20").Value = 202403.03
$ws.Range("N20").Value = 3968.69
$ws.Range("S20").Value = 202403.03
$ws.Range("H23").Value = 124882.54
$ws.Range("L23").Value = 117.87
$ws.Range("M23").Value = -18931.47
$ws.Range("H24").Value = 369255
$ws.Range("L24").Value = 131.76
$ws.Range("M24").Value = -89006.41
$ws.Range("H25").Value = 352204.5
$ws.Range("L25").Value = 123.21
$ws.Range("M25").Value = -66341.27
$ws.Range("H26").Value = 264835.89
$ws.Range("L26").Value = 119.74
$ws.Range("M26").Value = -43664.83
$ws.Range("H27").Value = 403946.75
$ws.Range("L27").Value = 111.59
$ws.Range("M27").Value = 173862.13
$ws.Range("N27").Value = 5608.46
$ws.Range("O27").Value = "Overpacing"
$ws.Range("S27").Value = 173862.13
$ws.Range("T27").Value = "NO"
$ws.Range("H28").Value = 246887.79
$ws.Range("L28").Value = 125.47
$ws.Range("M28").Value = -50110.43
$ws.Range("H29").Value = 198094.87
$ws.Range("L29").Value = 116.03
$ws.Range("M29").Value = -27364.86
$ws.Range("H30").Value = 426710.91
$ws.Range("L30").Value = 119.98
$ws.Range("M30").Value = -71050.64999999999
$ws.Range("H31").Value = 495019.74
$ws.Range("L31").Value = 118.18
$ws.Range("M31").Value = -76162.35000000001
$ws.Range("H32").Value = 337579.23
$ws.Range("L32").Value = 110.13
$ws.Range("M32").Value = 237161.04
$ws.Range("N32").Value = 4235.02
$ws.Range("O32").Value = "Overpacing"
$ws.Range("T32").Value = "NO"
$ws.Range("H33").Value = 85708.21000000001
$ws.Range("L33").Value = 128.62
$ws.Range("M33").Value = 43648.76
$ws.Range("N33").Value = 1364.02
$ws.Range("S33").Value = 43648.75999999999
$ws.Range("H34").Value = 700592.98
$ws.Range("L34").Value = 118.35
$ws.Range("M34").Value = -108642.49
$ws.Range("H35").Value = 480350.28
$ws.Range("L35").Value = 123.05
$ws.Range("M35").Value = -89968.47
$ws.Range("H36").Value = 150204.07
$ws.Range("L36").Value = 117.1
$ws.Range("M36").Value = -21935.28
$ws.Range("H38").Value = 436987.23
$ws.Range("L38").Value = 111.71
$ws.Range("M38").Value = -45790.66
$ws.Range("O38").Value = "Overpacing"
$ws.Range("H39").Value = 775972.6899999999
$ws.Range("L39").Value = 130.11
$ws.Range("M39").Value = -179593.76
$ws.Range("H40").Value = 71058.67
$ws.Range("L40").Value = 98.62
$ws.Range("M40").Value = 49030.78
$ws.Range("N40").Value = 1885.8
$ws.Range("S40").Value = 49030.78
$ws.Range("H43").Value = 450187.93
$ws.Range("L43").Value = 127.8
$ws.Range("M43").Value = -97930.47
$ws.Range("H45").Value = 305824.04
$ws.Range("L45").Value = 117.47
$ws.Range("M45").Value = 48786.97
$ws.Range("N45").Value = 2323.19
$ws.Range("S45").Value = 48786.97000000003

# --- Sheet: Feature_Importance (re-sorted by new Importance desc) ---
$ws2 = $wb.Worksheets.Item("Feature_Importance")
$ws2.Range("A2").Value = "Flight_Days"
$ws2.Range("B2").Value = 0.2612529686759745
$ws2.Range("A3").Value = "Days_Elapsed"
$ws2.Range("B3").Value = 0.2201912516487678
$ws2.Range("A4").Value = "Spend_to_Date"
$ws2.Range("B4").Value = 0.1277462596841187
$ws2.Range("A5").Value = "Total_Budget"
$ws2.Range("B5").Value = 0.1207321921315315
$ws2.Range("A6").Value = "Spend_Velocity"
$ws2.Range("B6").Value = 0.09635746771991167
$ws2.Range("A7").Value = "Pace_Ratio"
$ws2.Range("B7").Value = 0.08994689671939443
$ws2.Range("A8").Value = "DSP_enc"
$ws2.Range("B8").Value = 0.08377296342030142

# --- Sheet: Exec_Summary ---
$ws3 = $wb.Worksheets.Item("Exec_Summary")
$ws3.Range("B2").Value = 0.833
$ws3.Range("B3").Value = 1383187.49
$ws3.Range("B4").Value = "2026-02-09 12:45 UTC"
